$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H16").Value = -1
$ws.Range("H18").Value = -1
$ws.Range("H27").Value = -1
$ws.Range("H29").Value = -1
